$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 374.5
$ws.Range("I6").Value = 232.66667
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 698.00001
$ws.Range("L6").Value = 2400
$ws.Range("M6").Value = -586.00001
$ws.Range("N6").Value = -2624
$ws.Range("H17").Value = 1279.5652
$ws.Range("J17").Value = 1279.5652
$ws.Range("L17").Value = 3838.6956
$ws.Range("N17").Value = -4174.6956
$ws.Range("H100").Value = 50002090
$ws.Range("J100").Value = 125002220
$ws.Range("L100").Value = 125002220
$ws.Range("N100").Value = -125003302
$ws.Range("H103").Value = 12500576
$ws.Range("I103").Value = 429.75
$ws.Range("J103").Value = 16667292
$ws.Range("K103").Value = 1289.25
$ws.Range("L103").Value = 50001876
$ws.Range("M103").Value = -703.25
$ws.Range("N103").Value = -50003048
$ws.Range("H112").Value = 2962.5264
$ws.Range("J112").Value = 3202.303
$ws.Range("L112").Value = 9606.909
$ws.Range("N112").Value = -11822.909
$ws.Range("H118").Value = 400931.8
$ws.Range("I118").Value = 445036.66
$ws.Range("J118").Value = 3988
$ws.Range("K118").Value = 1335109.98
$ws.Range("L118").Value = 11964
$ws.Range("M118").Value = -1333452.98
$ws.Range("N118").Value = -15278
$ws.Range("H121").Value = 1045.7646
$ws.Range("J121").Value = 1045.7646
$ws.Range("L121").Value = 3137.2938
$ws.Range("N121").Value = -6631.293799999999
$ws.Range("H127").Value = 1142.2
$ws.Range("I127").Value = 388.5
$ws.Range("J127").Value = 1644.6666
$ws.Range("K127").Value = 1165.5
$ws.Range("L127").Value = 4933.9998
$ws.Range("M127").Value = 3794.5
$ws.Range("N127").Value = -14853.9998
$ws.Range("H135").Value = 304579.4
$ws.Range("I135").Value = 304579.4
$ws.Range("K135").Value = 2741214.6
$ws.Range("M135").Value = -2738679.6
$ws.Range("H137").Value = 2470.3635
$ws.Range("I137").Value = 1760.5714
$ws.Range("K137").Value = 5281.7142
$ws.Range("M137").Value = -2731.7142
$ws.Range("H138").Value = 2290.5945
$ws.Range("I138").Value = 2297.7693
$ws.Range("J138").Value = 2286.7083
$ws.Range("K138").Value = 6893.3079
$ws.Range("L138").Value = 6860.124899999999
$ws.Range("M138").Value = -1753.3079
$ws.Range("N138").Value = -17140.1249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2489.3333
$ws.Range("I61").Value = 2405
$ws.Range("K61").Value = 2405
$ws.Range("M61").Value = -2193
$ws.Range("H63").Value = 3353.1765
$ws.Range("J63").Value = 4136.909
$ws.Range("L63").Value = 4136.909
$ws.Range("N63").Value = -5508.909
$ws.Range("H66").Value = 3353.1765
$ws.Range("J66").Value = 4136.909
$ws.Range("L66").Value = 20684.545
$ws.Range("N66").Value = -27548.545
$ws.Range("H102").Value = 111113144
$ws.Range("I102").Value = 2184
$ws.Range("K102").Value = 2184
$ws.Range("M102").Value = -562
$ws.Range("H108").Value = 26407
$ws.Range("I108").Value = 23621
$ws.Range("J108").Value = 27800
$ws.Range("K108").Value = 23621
$ws.Range("L108").Value = 27800
$ws.Range("M108").Value = -19781
$ws.Range("N108").Value = -35480
$ws.Range("H110").Value = 11197.625
$ws.Range("I110").Value = 11811.2
$ws.Range("J110").Value = 10175
$ws.Range("K110").Value = 11811.2
$ws.Range("L110").Value = 10175
$ws.Range("M110").Value = -9766.200000000001
$ws.Range("N110").Value = -14265
$ws.Range("H123").Value = 2500000
$ws.Range("J123").Value = 2500000
$ws.Range("L123").Value = 2500000
$ws.Range("N123").Value = -2509800
$ws.Range("H128").Value = 56480
$ws.Range("J128").Value = 56480
$ws.Range("L128").Value = 56480
$ws.Range("N128").Value = -66440
$ws.Range("H136").Value = 2489.3333
$ws.Range("I136").Value = 2405
$ws.Range("K136").Value = 7215
$ws.Range("M136").Value = -4665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1850.3846
$ws.Range("I99").Value = 1178
$ws.Range("J99").Value = 4091.6667
$ws.Range("K99").Value = 1178
$ws.Range("L99").Value = 4091.6667
$ws.Range("M99").Value = 320
$ws.Range("N99").Value = -7087.6667
$ws.Range("H105").Value = 2660.1936
$ws.Range("I105").Value = 1922.091
$ws.Range("J105").Value = 4464.4443
$ws.Range("K105").Value = 1922.091
$ws.Range("L105").Value = 4464.4443
$ws.Range("M105").Value = -175.0909999999999
$ws.Range("N105").Value = -7958.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2424.919
$ws.Range("I31").Value = 1611.8889
$ws.Range("K31").Value = 1611.8889
$ws.Range("M31").Value = -1316.8889
$ws.Range("H34").Value = 2424.919
$ws.Range("I34").Value = 1611.8889
$ws.Range("K34").Value = 1611.8889
$ws.Range("M34").Value = -1409.8889
$ws.Range("H64").Value = 28000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 28000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 28000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -28496
$ws.Range("H67").Value = 28000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 28000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 28000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -29716
$ws.Range("H107").Value = 1181.0454
$ws.Range("I107").Value = 1344.8823
$ws.Range("J107").Value = 624
$ws.Range("K107").Value = 1344.8823
$ws.Range("L107").Value = 624
$ws.Range("M107").Value = 575.1177
$ws.Range("N107").Value = -4464
$ws.Range("H134").Value = 13085.538
$ws.Range("I134").Value = 14059.333
$ws.Range("J134").Value = 1400
$ws.Range("K134").Value = 42177.999
$ws.Range("L134").Value = 4200
$ws.Range("M134").Value = -39642.999
$ws.Range("N134").Value = -9270

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1972
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 3864.5715
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 11593.7145
$ws.Range("M25").Value = -1331
$ws.Range("N25").Value = -11931.7145
$ws.Range("H30").Value = 1972
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 3864.5715
$ws.Range("K30").Value = 1500
$ws.Range("L30").Value = 11593.7145
$ws.Range("M30").Value = -1398
$ws.Range("N30").Value = -11797.7145
$ws.Range("H60").Value = 1930.7693
$ws.Range("I60").Value = 100
$ws.Range("J60").Value = 4066.6667
$ws.Range("K60").Value = 300
$ws.Range("L60").Value = 12200.0001
$ws.Range("M60").Value = -49
$ws.Range("N60").Value = -12702.0001
$ws.Range("H68").Value = 3042.7144
$ws.Range("I68").Value = 7650
$ws.Range("K68").Value = 22950
$ws.Range("M68").Value = -22139
$ws.Range("H71").Value = 3042.7144
$ws.Range("I71").Value = 7650
$ws.Range("K71").Value = 68850
$ws.Range("M71").Value = -64794
$ws.Range("H121").Value = 742152.4399999999
$ws.Range("I121").Value = 640.5
$ws.Range("J121").Value = 954013
$ws.Range("K121").Value = 1921.5
$ws.Range("L121").Value = 2862039
$ws.Range("M121").Value = -611.5
$ws.Range("N121").Value = -2864659
$ws.Range("H131").Value = 887.20636
$ws.Range("J131").Value = 980.48
$ws.Range("L131").Value = 2941.44
$ws.Range("N131").Value = -13021.44

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 21333.334
$ws.Range("J47").Value = 21333.334
$ws.Range("L47").Value = 21333.334
$ws.Range("N47").Value = -22469.334
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7145232
$ws.Range("I7").Value = 9092659
$ws.Range("J7").Value = 4666.6665
$ws.Range("K7").Value = 9092659
$ws.Range("L7").Value = 4666.6665
$ws.Range("M7").Value = -9092547
$ws.Range("N7").Value = -4890.6665
$ws.Range("H95").Value = 35116
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 35116
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 35116
$ws.Range("M95").ClearContents()
$ws.Range("N95").Value = -40608
$ws.Range("H126").Value = 7145232
$ws.Range("I126").Value = 9092659
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 27277977
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -27275507
$ws.Range("N126").Value = -18939.9995
$ws.Range("H136").Value = 2443.125
$ws.Range("I136").Value = 1951.75
$ws.Range("J136").Value = 4900
$ws.Range("K136").Value = 5855.25
$ws.Range("L136").Value = 14700
$ws.Range("M136").Value = -3305.25
$ws.Range("N136").Value = -19800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1572.7084
$ws.Range("I126").Value = 712.5
$ws.Range("J126").Value = 2777
$ws.Range("K126").Value = 2137.5
$ws.Range("L126").Value = 8331
$ws.Range("M126").Value = 332.5
$ws.Range("N126").Value = -13271
$ws.Range("H128").Value = 47713.75
$ws.Range("J128").Value = 47713.75
$ws.Range("L128").Value = 47713.75
$ws.Range("N128").Value = -57673.75
$ws.Range("H136").Value = 2600.4614
$ws.Range("I136").Value = 2950.1
$ws.Range("J136").Value = 1435
$ws.Range("K136").Value = 8850.299999999999
$ws.Range("L136").Value = 4305
$ws.Range("M136").Value = -6300.299999999999
